$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TRUnits")

$newUnits = @(
    "SD_6_Droid_Company",
    "Storm_Cloud_Car_Company",
    "Mandalorian_Soldier_Company",
    "Mandalorian_Commando_Company",
    "Canderous_Assault_Tank_Company",
    "Mandalorian_Soldier_Team",
    "Mandalorian_Commando_Team",
    "Protocol_Legionnaire_Team",
    "Overracer_Speeder_Bike_Company",
    "SD_6_Hulk_Infantry_Droid_Company",
    "Police_Responder_Team",
    "Security_Trooper_Team",
    "PDF_Soldier_Team",
    "Military_Soldier_Team",
    "Scavenger_Team",
    "Heavy_Scavenger_Team",
    "Light_Mercenary_Team",
    "Mercenary_Team",
    "Elite_Mercenary_Team",
    "Navy_Trooper_Squad",
    "Enforcer_Trooper_Squad",
    "Army_Guard_Squad",
    "Army_Trooper_Squad",
    "Imperial_Medium_Armytrooper_Squad",
    "Imperial_Medic_Armytrooper_Squad",
    "Pentastar_Army_Trooper_Squad",
    "Pentastar_Medium_Armytrooper_Squad",
    "Pentastar_Medic_Armytrooper_Squad",
    "EVO_Trooper_Platoon",
    "Imperial_Compforce_Assault_Squad",
    "Imperial_Navy_Commando_Squad",
    "Army_Special_Missions_Squad",
    "Imperial_ISB_Infiltrator_Squad",
    "Imperial_IntSec_Operator_Squad",
    "Imperial_Stormtrooper_Squad",
    "Imperial_Medium_Stormtrooper_Squad",
    "Imperial_Medic_Stormtrooper_Squad",
    "Ysalamiri_Stormtrooper_Squad",
    "Imperial_Dark_Stormtrooper_Squad",
    "Imperial_Shadow_Stormtrooper_Squad",
    "Imperial_Royal_Guard_Stormtrooper_Squad",
    "Imperial_Elite_Guard_Stormtrooper_Squad",
    "Imperial_Nova_Stormtrooper_Squad",
    "Imperial_Incinerator_Stormtrooper_Squad",
    "Imperial_Army_Commando_Squad",
    "Imperial_Fleet_Commando_Squad",
    "Imperial_Fleet_Commando_Sniper_Squad",
    "Imperial_Guardian_Protector_Squad",
    "Imperial_Scout_Trooper_Assault_Squad",
    "Imperial_Scout_Trooper_Sniper_Squad",
    "Imperial_Scout_Trooper_Commando_Squad",
    "Zsinj_Raptor_Squad",
    "Zsinj_Raptor_Medic_Squad",
    "Storm_Commando_Platoon",
    "Imperial_Jumptrooper_Squad",
    "Imperial_EWEB_Company",
    "Imperial_74Z_Bike_Company",
    "Zsinj_74Z_Bike_Company",
    "Imperial_BARC_Company",
    "Imperial_64Y_Company",
    "Imperial_Dwarf_Spider_Droid_Squad",
    "Noghri_Assassin_Squad",
    "Imperial_Dark_Jedi_Squad",
    "Dark_Jedi_Squad",
    "Night_Sister_Company",
    "Imperial_Chrysalide_Company",
    "Imperial_PX10_Company",
    "Imperial_AT_PT_Company",
    "Imperial_AT_RT_Company",
    "Imperial_AT_MP_Company",
    "Imperial_AT_ST_Company",
    "Imperial_1L_Tank_Company",
    "Imperial_ULAV_Company",
    "Imperial_AT_ST_A_Company",
    "Imperial_ISP_Company",
    "Imperial_Repulsor_Scout_Company",
    "Imperial_RTT_Company",
    "Imperial_Gaba18_Group",
    "Imperial_INT4_Group",
    "Imperial_AT_DP_Company",
    "Imperial_Skiff_Company",
    "Imperial_Chariot_LAV_Company",
    "Imperial_TIE_Mauler_Company",
    "Imperial_Light_Artillery_Company",
    "Imperial_UT_AA_Company",
    "Imperial_AT_AP_Walker_Company",
    "Imperial_AT_AA_Company",
    "Imperial_Heavy_Artillery_Company",
    "Imperial_Missile_Artillery_Company",
    "Imperial_APC_Company",
    "Imperial_TX130_Company",
    "Imperial_1M_Tank_Company",
    "Imperial_Flashblind_Group",
    "Imperial_Century_Tank_Company",
    "Imperial_TNT_Company",
    "Imperial_2M_Company",
    "Imperial_S_1_Firehawke_Company",
    "Imperial_SA5_Company",
    "Imperial_Deathhawk_Group",
    "Imperial_SPMAG_Company",
    "Imperial_1H_Tank_Company",
    "Imperial_1H_Tank_Company_Green",
    "Imperial_SPMAT_Company",
    "Imperial_LAAT_Group",
    "Imperial_IDT_Group",
    "Imperial_MAAT_Group",
    "Imperial_Nemesis_Gunship_Group",
    "Imperial_B5_Juggernaut_Company",
    "Imperial_AT_TE_Walker_Company",
    "Imperial_PX4_Company",
    "Imperial_TRMB_Company",
    "Imperial_A5_Juggernaut_Company",
    "Imperial_A9_Fortress_Company",
    "Imperial_AT_AT_Company",
    "Imperial_C10_Siege_Tower_Company",
    "Imperial_Lancet_Group",
    "Imperial_Heavy_Recovery_Vehicle_Company",
    "Imperial_AT_AT_Refit_Company",
    "Imperial_XR85_Company",
    "Imperial_A6_Juggernaut_Company",
    "X34_Technical_Company",
    "ULAV_Early_Company",
    "Espo_Walker_Early_Squad",
    "ISP_Company",
    "AAC_1_Company",
    "Arrow_23_Company",
    "AA70_Company",
    "T1A_Company",
    "T2A_Company",
    "Freerunner_Early_Company",
    "Freerunner_Assault_Company",
    "Freerunner_AA_Company",
    "SA5_Company_Early_Company",
    "T3A_Company",
    "T4A_Company",
    "MZ8_Tank_Company",
    "A4_Juggernaut_Company",
    "Aratech_Battle_Platform_Company",
    "Skyhopper_Group",
    "Skyhopper_Antivehicle_Group",
    "Skyhopper_Primitive_Group",
    "Skyhopper_Security_Group",
    "Storm_Cloud_Car_Group",
    "JX30_Group",
    "Rebel_Infantry_Squad",
    "Rebel_Infantry_Squad_Deployed",
    "Rebel_Support_Infantry_Squad",
    "Defense_Trooper_Squad",
    "Rebel_Marine_Squad",
    "Rebel_Infiltrator_Team",
    "Jedi_Squad",
    "Rebel_AAC_2_Company",
    "Rebel_T1B_Company",
    "Rebel_AA5_Company",
    "Rebel_Snowspeeder_Wing",
    "Rebel_T2B_Company",
    "Rebel_Tracker_Company",
    "Rebel_MPTL_Company",
    "Rebel_Vwing_Group",
    "Rebel_T3B_Company",
    "Rebel_Gallofree_HTT_Company",
    "Rebel_Freerunner_Company",
    "Rebel_T4B_Company",
    "MAS_2xB_Company",
    "Phalanx_Commando_Team",
    "Phalanx_Trooper_Squad",
    "Phalanx_Medic_Squad",
    "Ysalamiri_Phalanx_Squad",
    "EotH_Scout_Brigade",
    "EotH_Kirov_Brigade",
    "Flame_Tank_Company",
    "AirStraeker_Company",
    "RFT_Brigade",
    "MMT_Brigade",
    "Mortar_Company",
    "Gilzean_Brigade",
    "CSA_B1_Droid_Squad",
    "CSA_Medium_Droid_Squad",
    "Espo_Squad",
    "CSA_Medium_Espo_Squad",
    "CSA_Medic_Espo_Squad",
    "CSA_Destroyer_Droid_Company",
    "X10_Group",
    "Class_I_Company",
    "Espo_Walker_Squad",
    "GX12_Company",
    "SX20_Company",
    "Hailfire_Company",
    "Persuader_Company",
    "JX40_Group",
    "Strikebreaker_Group",
    "K222_Group",
    "HRG_Commando_Squad",
    "Hapan_Infantry_Squad",
    "Maluri_Infantry_Squad",
    "Requud_Infantry_Squad",
    "Galney_Infantry_Squad",
    "Hapan_LightTank_Company",
    "Hapan_Artillery_Company",
    "Hapan_Transport_Company",
    "Hapan_Silanus_Company",
    "Hapan_HeavyTank_Company",
    "Hapan_Gunboat_Company",
    "SsiRuuk_Squad",
    "Yevetha_Infantry_Squad",
    "Yevethan_Plex_Squad",
    "Talon_Cloud_Car_Group"
)

$startRow = 276
for ($i = 0; $i -lt $newUnits.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newUnits[$i]
}

$ws.Range("A480").Select()
